$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "26.274.75"
$ws.Cells.Item(2,5).Value = "  -0.10%  "
$ws.Cells.Item(3,4).Value = "1.688.81"
$ws.Cells.Item(3,5).Value = "  +0.54%  "
$ws.Cells.Item(4,5).Value = "  +0.11%  "
$ws.Cells.Item(5,4).Value = "'217.63"
$ws.Cells.Item(5,5).Value = "  -0.30%  "
$ws.Cells.Item(6,4).Value = "'0.5383"
$ws.Cells.Item(6,5).Value = "  +2.52%  "
$ws.Cells.Item(7,5).Value = "  +0.08%  "
$ws.Cells.Item(8,4).Value = "'0.2727"
$ws.Cells.Item(8,5).Value = "  +1.24%  "
$ws.Cells.Item(9,4).Value = "'0.06425"
$ws.Cells.Item(9,5).Value = "  -0.66%  "
$ws.Cells.Item(10,4).Value = "'21.57"
$ws.Cells.Item(10,5).Value = "  -1.82%  "
$ws.Cells.Item(11,4).Value = "'0.07673"
$ws.Cells.Item(11,5).Value = "  +2.16%  "
$ws.Cells.Item(12,4).Value = "1.701.05"
$ws.Cells.Item(12,5).Value = "  +1.15%  "
$ws.Cells.Item(13,4).Value = "'4.525"
$ws.Cells.Item(13,5).Value = "  -0.07%  "
$ws.Cells.Item(14,4).Value = "'0.5774"
$ws.Cells.Item(14,5).Value = "  -0.57%  "
$ws.Cells.Item(15,4).Value = "'0.000008364"
$ws.Cells.Item(15,5).Value = "  -1.78%  "
$ws.Cells.Item(16,4).Value = "'66.71"
$ws.Cells.Item(16,5).Value = "  +2.93%  "
$ws.Cells.Item(17,4).Value = "26.334.49"
$ws.Cells.Item(17,5).Value = "  +0.02%  "
$ws.Cells.Item(18,4).Value = "'4.901"
$ws.Cells.Item(19,5).Value = "  +0.10%  "
$ws.Cells.Item(20,4).Value = "'10.85"
$ws.Cells.Item(20,5).Value = "  -0.26%  "
$ws.Cells.Item(21,4).Value = "'190.78"
$ws.Cells.Item(21,5).Value = "  +0.46%  "
$ws.Cells.Item(22,4).Value = "'6.256"
$ws.Cells.Item(22,5).Value = "  +0.83%  "
$ws.Cells.Item(23,5).Value = "  +0.04%  "
$ws.Cells.Item(24,4).Value = "'148.99"
$ws.Cells.Item(24,5).Value = "  +2.46%  "
$ws.Cells.Item(25,4).Value = "'0.1291"
$ws.Cells.Item(25,5).Value = "  +3.02%  "
$ws.Cells.Item(26,4).Value = "'7.843"
$ws.Cells.Item(26,5).Value = "  +0.21%  "
$ws.Cells.Item(27,4).Value = "'15.85"
$ws.Cells.Item(27,5).Value = "  +0.29%  "
$ws.Cells.Item(28,4).Value = "'0.06241"
$ws.Cells.Item(28,5).Value = "  -3.28%  "
$ws.Cells.Item(29,4).Value = "'1.368"
$ws.Cells.Item(29,5).Value = "  +0.55%  "
$ws.Cells.Item(30,4).Value = "'1.325"
$ws.Cells.Item(30,5).Value = "  +0.13%  "
$ws.Cells.Item(31,5).Value = "  -0.34%  "
$ws.Cells.Item(32,4).Value = "'3.579"
$ws.Cells.Item(32,5).Value = "  -0.41%  "
$ws.Cells.Item(33,4).Value = "'1.671"
$ws.Cells.Item(33,5).Value = "  +0.32%  "
$ws.Cells.Item(34,4).Value = "'1.029"
$ws.Cells.Item(34,5).Value = "  -0.04%  "
$ws.Cells.Item(35,4).Value = "'0.6177"
$ws.Cells.Item(35,5).Value = "  -0.97%  "
$ws.Cells.Item(36,5).Value = "  +0.49%  "
$ws.Cells.Item(37,4).Value = "'2.759"
$ws.Cells.Item(37,5).Value = "  +0.81%  "
$ws.Cells.Item(38,4).Value = "'0.01650"
$ws.Cells.Item(38,5).Value = "  +1.53%  "
$ws.Cells.Item(39,2).Value = "FraxShare"
$ws.Cells.Item(39,3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(39,4).Value = "'6.113"
$ws.Cells.Item(39,5).Value = "  -5.02%  "
$ws.Cells.Item(40,2).Value = "Maker"
$ws.Cells.Item(40,3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(40,4).Value = "1.106.26"
$ws.Cells.Item(40,5).Value = "  -0.13%  "
$ws.Cells.Item(41,4).Value = "'0.8786"
$ws.Cells.Item(41,5).Value = "  +0.34%  "
$ws.Cells.Item(42,5).Value = "  -0.20%  "
$ws.Cells.Item(43,4).Value = "'101.08"
$ws.Cells.Item(43,5).Value = "  +0.38%  "
$ws.Cells.Item(44,4).Value = "1.841.16"
$ws.Cells.Item(44,5).Value = "  +0.62%  "
$ws.Cells.Item(45,5).Value = "  +0.37%  "
$ws.Cells.Item(46,4).Value = "'57.57"
$ws.Cells.Item(46,5).Value = "  +1.09%  "
$ws.Cells.Item(47,4).Value = "'8.124"
$ws.Cells.Item(47,5).Value = "  -0.77%  "
$ws.Cells.Item(48,5).Value = "  -0.47%  "
$ws.Cells.Item(49,4).Value = "'0.05286"
$ws.Cells.Item(49,5).Value = "  +0.35%  "
$ws.Cells.Item(50,5).Value = "  +0.06%  "
$ws.Cells.Item(51,4).Value = "'6.041"
$ws.Cells.Item(51,5).Value = "  -0.70%  "
